# Add a new "Jira" column (F) to the KeywordFramework sheet, with the
# Jira ticket numbers lined up against each TestCase header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1: reuse the existing header style (same as A1:E1) ---
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data cells F2:F9: reuse the existing body style (same as B2:B9) ---
$ws.Range("B2:B9").Copy()
$ws.Range("F2:F9").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Values ---
$ws.Range("F1").Value = "Jira"
$ws.Range("F2").Value = "Qa-0124"
$ws.Range("F6").Value = "Qa-0125"

# --- Update the active selection to reflect where the user ended up ---
[void]$ws.Range("M17").Select()
